$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata" ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: fhir/fr/medication -> ig/fhir/medication
$meta.Range("B2").Value = "https://hl7.fr/ig/fhir/medication/ValueSet/fr-medication-reconciliation-document-type"

# Name: FrMedicationReconciliationDocumentType -> FRMedicationReconciliationDocumentType
$meta.Range("B4").Value = "FRMedicationReconciliationDocumentType"

# Title: InterOp'Santé -> Interop'Santé
$meta.Range("B5").Value = "value set Interop'Santé - type de document de la ressource Composition d'une FCT"

# Date: updated timestamp
$meta.Range("B8").Value = "2026-01-15T08:54:26+00:00"

# Jurisdiction: (empty) -> FRANCE
$meta.Range("B11").Value = "FRANCE"

# --- Sheet "Include #0" ---
$inc = $wb.Worksheets.Item("Include #0")

# System URI value: fhir/fr/medication -> ig/fhir/medication
$inc.Range("B4").Value = "https://hl7.fr/ig/fhir/medication/CodeSystem/fr-document-type"
